# "correction bus group debit calcul"
#
# 1. On sheet "map s2": clear G7 and change C8 from 0 to "P", and leave the
#    selection there at C8.
# 2. Add a new sheet "map T" (copied from "map (3)" so it inherits the same
#    column widths/styles/conditional formatting), positioned after
#    "map (3)", then tweak the handful of cells that differ from its
#    source, and make it the active sheet/selection (D9).

$wb = $excel.ActiveWorkbook

# --- 1. Fix up "map s2" --------------------------------------------------
$mapS2 = $wb.Worksheets.Item("map s2")
$mapS2.Range("G7").Clear()
$mapS2.Range("C8").Value = "P"
$mapS2.Range("C8").Select()

# --- 2. Duplicate "map (3)" into a new "map T" sheet at the end ---------
$src = $wb.Worksheets.Item("map (3)")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "map T"

# Cells that differ from the "map (3)" source data
$newSheet.Range("C2").Value = "C"
$newSheet.Range("E2").Value = 1
$newSheet.Range("G2").Value = "C"
$newSheet.Range("D4").Value = "E"
$newSheet.Range("C6").Value = 0
$newSheet.Range("C7").Value = "T"
$newSheet.Range("E7").Value = "P"
$newSheet.Range("F7").Value = 0
$newSheet.Range("G7").Value = "T"
$newSheet.Range("C9").Value = "E"
$newSheet.Range("D9").Value = 0
$newSheet.Range("G9").Value = "E"
$newSheet.Range("C11").Value = "C"
$newSheet.Range("E11").Value = 1
$newSheet.Range("G11").Value = "C"

# --- 3. View state: "map T" becomes the active/selected sheet -----------
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 205
$newSheet.Range("D9").Select()
